$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows correspond to sheet rows 2-25 (A column 0-23)
# Columns correspond to sheet columns B-N
$rows = @(
    @(1.438419810470521, 0.2238222309472064, 0.07798299849743273, 0.06710533727783652, 0, 1.373402489101181, 1.233393531797105, 1.11533496149174, 0, 0, 0.2138799314047475, 0.3024060800092414, 1.636194416556162),
    @(1.336932980063239, 0.1952902509472096, 0.07080067220292108, 0.06695115771748128, 0, 1.351429827076004, 1.229531673347481, 1.114865641306707, 0, 0, 0.2110881171394254, 0.2869623549792806, 1.65664991341413),
    @(1.275266687061389, 0.1777265266140091, 0.06643052715951114, 0.06685702267133298, 0, 1.338928171771158, 1.227862139483236, 1.115244843096036, 0, 0, 0.2094845302230723, 0.2776347045117831, 1.669831875395477),
    @(1.250299809336866, 0.1705575481012715, 0.06465955801908763, 0.06681879348147834, 0, 1.334081418448363, 1.227357874958784, 1.115566762276295, 0, 0, 0.2088588660681481, 0.2738725404491547, 1.675359880899757),
    @(1.246163899826513, 0.1693664340117493, 0.0643660845444316, 0.06681245344172293, 0, 1.333291552125317, 1.227284766988902, 1.115630314434064, 0, 0, 0.208756654716602, 0.27325018695187, 1.676287233005294),
    @(1.27492931658071, 0.1776298907064984, 0.06640660331395054, 0.06685650656860798, 0, 1.338861804974115, 1.227854626316244, 1.115248507416482, 0, 0, 0.2094759796929608, 0.2775838090419782, 1.669905795651303),
    @(1.403292782769654, 0.213993543438562, 0.07549819311945782, 0.06705206443997724, 0, 1.365620249904538, 1.231916113605678, 1.115034369050385, 0, 0, 0.2128943572797652, 0.2970489294964196, 1.643118209007953),
    @(1.660165049440423, 0.2849630605430207, 0.09364899676756977, 0.06743989094694047, 0, 1.425999531134863, 1.245466281647367, 1.119930406894106, 0, 0, 0.2204760189410564, 0.3364516739918173, 1.595533300262019),
    @(1.852076577010337, 0.3369263068685768, 0.1071907457053385, 0.06772766796182861, 0, 1.475261145179076, 1.258855170439972, 1.126799290534123, 0, 0, 0.2265835801247533, 0.3661593949402047, 1.563598130105291),
    @(1.940085867593666, 0.360533627172515, 0.1133982744091213, 0.06785924738129889, 0, 1.498753784633607, 1.265698117228226, 1.130641464957009, 0, 0, 0.2294792305696802, 0.3798408183545945, 1.549729264962149),
    @(1.97351497812258, 0.3694690710851205, 0.1157558601856294, 0.06790917232127569, 0, 1.507806954203431, 1.268398013982079, 1.132200100180626, 0, 0, 0.2305926266167546, 0.3850457497687785, 1.544572438882231),
    @(1.966310886730241, 0.3675448420942757, 0.1152478017196472, 0.06789841566195332, 0, 1.50585018888799, 1.267811704784947, 1.131859800441035, 0, 0, 0.2303520861215844, 0.3839237035555811, 1.545678823016758),
    @(1.94283405758199, 0.3612688325347335, 0.1135920948152034, 0.0678633527436423, 0, 1.499495440336773, 1.265918059853874, 1.130767613807876, 0, 0, 0.2295704920508683, 0.3802685485730493, 1.549303103819064),
    @(1.928467096550662, 0.3574240665552679, 0.1125788326654771, 0.06784188862155283, 0, 1.495623458642569, 1.264772306679419, 1.130112136428167, 0, 0, 0.2290939411870454, 0.3780327964681831, 1.551535463025186),
    @(1.846339215408818, 0.3353829164194906, 0.1067860319439262, 0.06771908260231818, 0, 1.473747750293001, 1.258423143363729, 1.12656267681411, 0, 0, 0.2263967035592884, 0.3652686489947499, 1.564517771843136),
    @(1.796137781795323, 0.3218536693285898, 0.1032445444585193, 0.06764391804309078, 0, 1.460606046909504, 1.254721127495742, 1.124569327118323, 0, 0, 0.2247720831862239, 0.3574811073146904, 1.572650897588108),
    @(1.767329855419291, 0.3140690647066435, 0.1012120208961846, 0.06760074811251737, 0, 1.453149175627004, 1.25266260649181, 1.123490307125678, 0, 0, 0.2238486837781721, 0.3530176582839601, 1.57739082867451),
    @(1.757587423424184, 0.3114328166453504, 0.1005246032208191, 0.06758614221348225, 0, 1.450641868990147, 1.251977769816591, 1.123136547464199, 0, 0, 0.2235379323233104, 0.3515091135895148, 1.579006322470978),
    @(1.801474915329607, 0.3232941826559852, 0.1036210807602487, 0.06765191292673478, 0, 1.461994449613996, 1.255107884702994, 1.124774532711434, 0, 0, 0.2249438842089262, 0.3583084753231347, 1.571778696407552),
    @(1.949727004705608, 0.3631123583144813, 0.1140782270059475, 0.06787364887387337, 0, 1.501357712680317, 1.266471317739928, 1.131085597334888, 0, 0, 0.2297996070802668, 0.3813415026330134, 1.548235983622778),
    @(2.047212350456164, 0.3891118842391279, 0.1209530250766875, 0.06801914251234642, 0, 1.527999679580631, 1.274531317543563, 1.135814853421635, 0, 0, 0.2330714746224487, 0.3965352490597809, 1.53340340974723),
    @(1.995128186443708, 0.3752375312863023, 0.1172800733514237, 0.06794143620613013, 0, 1.513696151443924, 1.270171444553654, 1.133235267324622, 0, 0, 0.2313162127703805, 0.3884132110739245, 1.541269044469711),
    @(1.799061830479616, 0.32264294673314, 0.1034508377306054, 0.06764829830594365, 0, 1.461366446125766, 1.254932814325002, 1.124681550652404, 0, 0, 0.2248661798991805, 0.3579343795366654, 1.572172818947781),
    @(1.590117982065522, 0.2657972971276763, 0.08870320584213687, 0.06733448916795748, 0, 1.408810582879454, 1.241199707189054, 1.118033602743687, 0, 0, 0.2183307628495896, 0.3256595459041307, 1.607875763147124)
)

$data = New-Object "object[,]" 24,13
for ($i = 0; $i -lt $rows.Count; $i++) {
    for ($j = 0; $j -lt $rows[$i].Count; $j++) {
        $data[$i,$j] = $rows[$i][$j]
    }
}

$ws.Range("B2:N25").Value = $data

